# Weekly crime data refresh (week of 12/11/2023 - 12/17/2023), Volume 30 Number 50
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ==== Header text changes ====
$ws.Range("A8").Value = "Volume 30   Number  50"
$ws.Range("C9").Value = "Report Covering the Week  12/11/2023  Through  12/17/2023"

# ==== Numeric / simple string cell updates (rows 14-30) ====
$ws.Range("E14").Value = "***.*"
$ws.Range("L14").Value = -66.666666666666
$ws.Range("N14").Value = -90.410958904109
$ws.Range("E15").Value = -100
$ws.Range("F15").Value = 3
$ws.Range("H15").Value = -25
$ws.Range("I15").Value = 35
$ws.Range("J15").Value = 35
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 6.060606060606
$ws.Range("M15").Value = 2.941176470588
$ws.Range("N15").Value = -58.333333333333
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -16.666666666666
$ws.Range("F16").Value = 18
$ws.Range("G16").Value = 28
$ws.Range("H16").Value = -35.714285714285
$ws.Range("I16").Value = 264
$ws.Range("J16").Value = 335
$ws.Range("K16").Value = -21.194029850746
$ws.Range("L16").Value = 20
$ws.Range("M16").Value = -45.341614906832
$ws.Range("N16").Value = -87.811634349030
$ws.Range("C17").Value = 11
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = -8.333333333333
$ws.Range("I17").Value = 656
$ws.Range("J17").Value = 699
$ws.Range("K17").Value = -6.151645207439
$ws.Range("L17").Value = 5.977382875605
$ws.Range("M17").Value = 17.985611510791
$ws.Range("N17").Value = -40.687160940325
$ws.Range("C18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -10
$ws.Range("I18").Value = 122
$ws.Range("J18").Value = 187
$ws.Range("K18").Value = -34.759358288770
$ws.Range("L18").Value = -15.277777777777
$ws.Range("M18").Value = -55.311355311355
$ws.Range("N18").Value = -83.776595744680
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -11.111111111111
$ws.Range("F19").Value = 30
$ws.Range("H19").Value = -16.666666666666
$ws.Range("I19").Value = 325
$ws.Range("J19").Value = 400
$ws.Range("K19").Value = -18.75
$ws.Range("L19").Value = 3.174603174603
$ws.Range("M19").Value = -25.799086757990
$ws.Range("N19").Value = -64.673913043478
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 145
$ws.Range("J20").Value = 183
$ws.Range("K20").Value = -20.765027322404
$ws.Range("L20").Value = -10.493827160493
$ws.Range("M20").Value = 2.112676056338
$ws.Range("N20").Value = -79.015918958031
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 33
$ws.Range("E21").Value = -21.212121212121
$ws.Range("F21").Value = 112
$ws.Range("G21").Value = 124
$ws.Range("H21").Value = -9.677419354838
$ws.Range("I21").Value = 1554
$ws.Range("J21").Value = 1860
$ws.Range("K21").Value = -16.451612903225
$ws.Range("L21").Value = 2.642007926023
$ws.Range("M21").Value = -20.389344262295
$ws.Range("N21").Value = -73.169889502762
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 8
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = 33.333333333333
$ws.Range("I22").Value = 51
$ws.Range("J22").Value = 48
$ws.Range("K22").Value = 6.25
$ws.Range("L22").Value = 18.604651162790
$ws.Range("M22").Value = -8.928571428571
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 10
$ws.Range("E23").Value = -60
$ws.Range("F23").Value = 25
$ws.Range("G23").Value = 21
$ws.Range("H23").Value = 19.047619047619
$ws.Range("I23").Value = 329
$ws.Range("J23").Value = 368
$ws.Range("K23").Value = -10.597826086956
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 34.836065573770
$ws.Range("C24").Value = 32
$ws.Range("D24").Value = 26
$ws.Range("E24").Value = 23.076923076923
$ws.Range("F24").Value = 98
$ws.Range("G24").Value = 91
$ws.Range("H24").Value = 7.692307692307
$ws.Range("I24").Value = 1051
$ws.Range("J24").Value = 1113
$ws.Range("K24").Value = -5.570530098831
$ws.Range("L24").Value = 33.885350318471
$ws.Range("M24").Value = 8.911917098445
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 19
$ws.Range("E25").Value = -26.315789473684
$ws.Range("F25").Value = 50
$ws.Range("G25").Value = 61
$ws.Range("H25").Value = -18.032786885245
$ws.Range("I25").Value = 916
$ws.Range("J25").Value = 911
$ws.Range("K25").Value = 0.548847420417
$ws.Range("L25").Value = 11.300121506682
$ws.Range("M25").Value = -26.484751203852
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 4
$ws.Range("G26").Value = 6
$ws.Range("H26").Value = -33.333333333333
$ws.Range("I26").Value = 43
$ws.Range("J26").Value = 58
$ws.Range("K26").Value = -25.862068965517
$ws.Range("L26").Value = -27.118644067796
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -14.285714285714
$ws.Range("I27").Value = 80
$ws.Range("J27").Value = 84
$ws.Range("K27").Value = -4.761904761904
$ws.Range("L27").Value = -8.045977011494
$ws.Range("C28").Value = 2
$ws.Range("E28").Value = "***.*"
$ws.Range("F28").Value = 3
$ws.Range("H28").Value = 200
$ws.Range("I28").Value = 46
$ws.Range("K28").Value = -38.666666666666
$ws.Range("L28").Value = -48.888888888888
$ws.Range("M28").Value = -54.455445544554
$ws.Range("N28").Value = -85.714285714285
$ws.Range("C29").Value = 2
$ws.Range("E29").Value = "***.*"
$ws.Range("F29").Value = 3
$ws.Range("H29").Value = 200
$ws.Range("I29").Value = 42
$ws.Range("K29").Value = -28.813559322033
$ws.Range("L29").Value = -45.454545454545
$ws.Range("M29").Value = -48.780487804878
$ws.Range("N29").Value = -85.714285714285
$ws.Range("H30").Value = "***.*"
$ws.Range("L30").Value = -66.666666666666

# ==== Text cells whose content looks numeric (must force text type) ====
$styleDonor = $ws.Range("C14")
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "0"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "0"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "0"
$styleDonor.Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("G30").PasteSpecial(-4122)
